$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 138, shifting existing row 138 (and below) down to 139.
$ws.Rows.Item(138).Insert()

# Populate the newly inserted row 138 with the new record's data.
$ws.Range("A138").Value = 6
$ws.Range("B138").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C138").Value = "Metropolitana"
$ws.Range("D138").Value = 44762
$ws.Range("D138").NumberFormat = $ws.Range("D139").NumberFormat
$ws.Range("E138").Value = 13
$ws.Range("F138").Value = 100112001
$ws.Range("G138").Value = "Berenjena"
$ws.Range("H138").Value = "Sin especificar"
$ws.Range("I138").Value = "Primera"
$ws.Range("J138").Value = 400
$ws.Range("K138").Value = 10000
$ws.Range("L138").Value = 11000
$ws.Range("M138").Value = 10425
$ws.Range("N138").Value = "`$/caja 50 unidades"
$ws.Range("O138").Value = "Región de Arica y Parinacota"
$ws.Range("P138").Value = 208
$ws.Range("Q138").Value = 50
$ws.Range("R138").Value = "Hortaliza"
